$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Record the new selection on Sheet1 before switching sheets (matches target diff)
$ws1.Range("I1").Select() | Out-Null

# Add Sheet2 right after Sheet1
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)

# --- Copy cell formatting (styles) from Sheet1 onto the same layout in Sheet2 ---
# Columns A-D line up directly between the two sheets
$ws1.Range("A1:D31").Copy()
$ws2.Range("A1").PasteSpecial(-4122)
# Sheet1's E,F,G (Live/Equiv./Invalid) become Sheet2's F,G,H (E is the new Trivial+Killed sum)
$ws1.Range("E1:G31").Copy()
$ws2.Range("F1").PasteSpecial(-4122)
# Give column E (new) the same look as the other data columns (copy style from D)
$ws1.Range("D1:D31").Copy()
$ws2.Range("E1").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Row 1 header values ---
$ws2.Range("A1").Value = "Mutation Operator"
$ws2.Range("B1").Value = "Gen."
$ws2.Range("C1").Value = "Trivial"
$ws2.Range("D1").Value = "Killed"
$ws2.Range("F1").Value = "Live"
$ws2.Range("G1").Value = "Equiv."
$ws2.Range("H1").Value = "Invalid"

# --- Data values, copied straight from Sheet1 (same operator names, Gen/Trivial/Killed counts) ---
$ws1.Range("A2:D31").Copy()
$ws2.Range("A2").PasteSpecial(-4163)
$ws1.Range("E2:G31").Copy()
$ws2.Range("F2").PasteSpecial(-4163)

$excel.CutCopyMode = 0

# --- New column E: Trivial + Killed ---
$ws2.Range("E2").Formula = "=C2+D2"
$ws2.Range("E3:E31").Formula = "=C3+D3"

# --- Column widths: wide operator name column, hide the raw helper columns C:D ---
$ws2.Columns("A:A").ColumnWidth = 28.25
$ws2.Columns("C:D").ColumnWidth = 0
$ws2.Columns("C:D").Hidden = $true

# Final selection/scroll position on the new sheet
$ws2.Range("I12").Select() | Out-Null
